$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 149, pushing the existing rows 149:165 down to 150:166.
$ws.Rows(149).Insert()

# Populate the newly inserted row 149 with the new weekly price record.
$ws.Range("A149").Value = 7
$ws.Range("B149").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C149").Value = "Ñuble"
$ws.Range("D149").Value = 45127
$ws.Range("E149").Value = 16
$ws.Range("F149").Value = "Fruta"
$ws.Range("G149").Value = 100108
$ws.Range("H149").Value = "Tropicales y subtropicales"
$ws.Range("I149").Value = 100108002
$ws.Range("J149").Value = "Mango"
$ws.Range("K149").Value = "Sin especificar"
$ws.Range("L149").Value = "Primera"
$ws.Range("M149").Value = 40
$ws.Range("N149").Value = 9000
$ws.Range("O149").Value = 9000
$ws.Range("P149").Value = 9000
$ws.Range("Q149").Value = "$/bandeja 4 kilos"
$ws.Range("R149").Value = "Brasil"
$ws.Range("S149").Value = 2250
$ws.Range("T149").Value = 4
